$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.921.44"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").Value = "2.902.44"
$ws.Range("E3").Value = "  -1.16%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.49"
$ws.Range("E5").Value = "  -3.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.64"
$ws.Range("E6").Value = "  -2.50%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("E8").Value = "  -0.95%  "

$ws.Range("D9").Value = "2.899.97"
$ws.Range("E9").Value = "  -1.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.96"
$ws.Range("E10").Value = "  -2.84%  "

$ws.Range("E11").Value = "  -2.73%  "

$ws.Range("E12").Value = "  -1.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.13"
$ws.Range("E14").Value = "  -0.45%  "

$ws.Range("E15").Value = "  -0.15%  "

$ws.Range("D16").Value = "3.382.71"

$ws.Range("D17").Value = "61.857.04"

$ws.Range("D18").Value = "2.902.53"
$ws.Range("E18").Value = "  -1.17%  "

$ws.Range("E19").Value = "  -1.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "430.04"

$ws.Range("E21").Value = "  -4.07%  "

$ws.Range("E22").Value = "  -1.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.87"

$ws.Range("E24").Value = "  -1.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.00"
$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("E26").Value = "  -8.17%  "

$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("E28").Value = "  -3.12%  "

$ws.Range("E29").Value = "  +7.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.99"
$ws.Range("E30").Value = "  -3.91%  "

$ws.Range("E31").Value = "  -2.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.03"
$ws.Range("E32").Value = "  -6.09%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("E34").Value = "  -3.01%  "

$ws.Range("E35").Value = "  -1.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.956"
$ws.Range("E36").Value = "  -3.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.83"
$ws.Range("E38").Value = "  -1.57%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.85"
$ws.Range("E39").Value = "  -5.91%  "

$ws.Range("E40").Value = "  -4.78%  "

$ws.Range("E41").Value = "  -1.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.96"
$ws.Range("E43").Value = "  +4.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.269"
$ws.Range("E44").Value = "  -2.09%  "

$ws.Range("D45").Value = "2.705.31"
$ws.Range("E45").Value = "  +0.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0337"
$ws.Range("E46").Value = "  -0.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "131.78"
$ws.Range("E47").Value = "  -2.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "346.29"
$ws.Range("E48").Value = "  -1.69%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("E50").Value = "  -0.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.58"
$ws.Range("E51").Value = "  -2.41%  "
